$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the current (pre-edit) values of row 16 and row 17 for the columns
# that differ between the two rows, then swap them.

$cols = @("A","B","C","D","E","F","G","H","I","J","K","P","Q","R")

$row16 = @{}
$row17 = @{}
foreach ($col in $cols) {
    $row16[$col] = $ws.Range("${col}16").Value2
    $row17[$col] = $ws.Range("${col}17").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}16").Value2 = $row17[$col]
    $ws.Range("${col}17").Value2 = $row16[$col]
}

# The empty "L" placeholder cell moves from row 16 to row 17.
$ws.Range("L16").Value2 = $null
$ws.Range("L17").Value2 = ""
